$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.962330102920532
$ws.Range("B1").Value = 2.577364206314087
$ws.Range("C1").Value = 2.338207006454468
$ws.Range("D1").Value = 2.465882062911987
$ws.Range("E1").Value = 3.186347723007202
